$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Adicionada introdução do experimento" - label a handful of helper cells
# to the right of the data table (columns E:G) that introduce/annotate the
# experiment data already in A:C.
#
# Order matters here: cells are written in the same order their text first
# appears (ref, -, V) so the shared-string table comes out in that order.

# E15: "ref" reference marker, styled like the red-font rows (style used by
# e.g. B15/C15 - red, non-bold font).
$ws.Range("E15").Value = "ref"
$ws.Range("E15").Font.Color = 255

# E7 / E24: "-" placeholder markers, default formatting.
$ws.Range("E7").Value = "-"
$ws.Range("E24").Value = "-"

# F6 / G6: "V" (Volt) markers, default formatting.
$ws.Range("F6").Value = "V"
$ws.Range("G6").Value = "V"

# F23: formatted like the bold red rows (style used by e.g. A23/B23/C23)
# but left without a value.
$ws.Range("F23").Font.Color = 255
$ws.Range("F23").Font.Bold = $true

# Leave the selection where the author ended up after these edits.
[void]$ws.Range("F7").Select()
